$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (column widths/positions unchanged, just renamed
# vessel length buckets per the "cleaning nvessels tables" commit).
$ws.Range("B1").Value = "Up to 24 feet"
$ws.Range("D1").Value = "40 to 61 feet"
$ws.Range("E1").Value = "65 to 84 feet"

# Update the view: scroll so column C is the leftmost visible column and
# move the active selection to H1.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H1").Select()
